# Update published master datasets (2026-02-18)
# Adds the new "Southwest Wales Net Zero Industry Launchpad Round 2 -
# Call for challenge holders" opportunity (Innovate UK Business Connect)
# to the "This weeks opportunities" and "All calls" sheets, rolls the new
# week forward on "Weekly counts (Wed)", and bumps the row-count totals on
# "Sense check".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. "This weeks opportunities" - insert the new opportunity as row 2
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("This weeks opportunities")

$ws1.Range("A2").Value = "Innovate UK Business Connect"
$ws1.Range("B2").Value = "Southwest Wales Net Zero Industry Launchpad Round 2 - Call for challenge holders"
$ws1.Range("C2").Value = "https://iuk-business-connect.org.uk/opportunities/southwest-wales-net-zero-industry-launchpad-round-2-call-for-challenge-holders/"
$ws1.Range("D2").Value = "2026-02-18 10:32"
$ws1.Range("E2").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Range("E2").Value = 46071.43888888889
$ws1.Range("F2").Value = $true
$ws1.Range("H2").Value = $false
$ws1.Range("I2").Value = "18/02/2026"
$ws1.Range("J2").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Range("J2").Value = 46071
$ws1.Range("K2").Value = "13/03/2026                              00:00"
$ws1.Range("L2").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Range("L2").Value = 46094
$ws1.Range("M2").Value = [char]0x00A3 + "25,000"
$ws1.Range("N2").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Range("N2").Value = 46071

# ---------------------------------------------------------------
# 2. "Weekly counts (Wed)" - insert a new week-commencing row at the top
#    of the data (row 2), pushing all the other weeks down by one row.
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Weekly counts (Wed)")

$ws2.Rows.Item(2).Insert()
$ws2.Range("A2").NumberFormat = "YYYY-MM-DD"
$ws2.Range("A2").Value = 46071
$ws2.Range("B2").Value = 0
$ws2.Range("C2").Value = 1
$ws2.Range("D2").Value = 1

# ---------------------------------------------------------------
# 3. "All calls" - insert the new opportunity as row 147 (in
#    Date-Added-descending order), pushing every later row down by one.
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("All calls")

$ws3.Rows.Item(147).Insert()
$ws3.Range("A147").Value = "Innovate UK Business Connect"
$ws3.Range("B147").Value = "Southwest Wales Net Zero Industry Launchpad Round 2 - Call for challenge holders"
$ws3.Range("C147").Value = "https://iuk-business-connect.org.uk/opportunities/southwest-wales-net-zero-industry-launchpad-round-2-call-for-challenge-holders/"
$ws3.Range("D147").Value = "2026-02-18 10:32"
$ws3.Range("E147").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws3.Range("E147").Value = 46071.43888888889
$ws3.Range("F147").Value = $true
$ws3.Range("H147").Value = $false
$ws3.Range("I147").Value = "18/02/2026"
$ws3.Range("J147").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws3.Range("J147").Value = 46071
$ws3.Range("K147").Value = "13/03/2026                              00:00"
$ws3.Range("L147").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws3.Range("L147").Value = 46094
$ws3.Range("M147").Value = [char]0x00A3 + "25,000"
$ws3.Range("N147").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws3.Range("N147").Value = 46071

# ---------------------------------------------------------------
# 4. "Sense check" - bump the raw-row-count totals now that one more
#    Business Connect row (and therefore one more combined row) exists.
# ---------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Sense check")

$ws4.Range("B2").Value = 146
$ws4.Range("C2").Value = 146
$ws4.Range("B4").Value = 250
$ws4.Range("C4").Value = 250
